# "Summing the Empty Values" - fill in the previously-empty Total columns
# (E = Share Capital + Reserves, K = Sales + Other Income) for rows 5-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> (value, format-source cell already carrying the shared
# "centered" style used throughout the data block rows 5-10). The format
# source is a neighboring, already-populated cell in the same row.
$totals = @(
    @{ Addr = "E5";  Value = 3.76;               FormatFrom = "D5" }
    @{ Addr = "K5";  Value = 13.37;              FormatFrom = "J5" }
    @{ Addr = "E6";  Value = 9.040000000000001;  FormatFrom = "D6" }
    @{ Addr = "K6";  Value = 21.89;              FormatFrom = "J6" }
    @{ Addr = "E7";  Value = 13.94;              FormatFrom = "D7" }
    @{ Addr = "K7";  Value = 36.38;              FormatFrom = "J7" }
    @{ Addr = "E8";  Value = 23.36;              FormatFrom = "D8" }
    @{ Addr = "K8";  Value = 21.11;              FormatFrom = "J8" }
    @{ Addr = "E9";  Value = 25.04;              FormatFrom = "D9" }
    @{ Addr = "K9";  Value = 22.7;               FormatFrom = "J9" }
    @{ Addr = "E10"; Value = 26.15;              FormatFrom = "D10" }
    @{ Addr = "K10"; Value = 31.89;              FormatFrom = "J10" }
)

foreach ($item in $totals) {
    $ws.Range($item.FormatFrom).Copy()
    $ws.Range($item.Addr).PasteSpecial(-4122)
    $ws.Range($item.Addr).Value = $item.Value
}
